# Helper: build "tok tok tok ..." (count repetitions of tok separated by single spaces)
# without using the "+"/"+=" operators on strings, since this interpreter auto-coerces
# decimal-looking string operands (e.g. "00", "07") to numeric addition.
function RepeatToken($token, $count) {
    $result = ""
    for ($i = 0; $i -lt $count; $i++) {
        if ($i -gt 0) {
            $result = "$result $token"
        } else {
            $result = "$token"
        }
    }
    return $result
}

$d = $word.ActiveDocument

# Locate the end of the existing "2 messages with length 6 and channel width of 20"
# example (its second data line), which is where the new test sample should be
# inserted, mirroring the surrounding "blank line" separated block structure used
# throughout this section of the document.
$anchor = $d.Content
$found = $anchor.Find.Execute(
    "2, 1, 07 07 07 07 07 07 01 D5 FF FF FF FF FF FF AA BB CC DD 00 00",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    Write-Host "ERROR: anchor paragraph not found"
} else {
    $anchor.Collapse(0)

    # Insert four new empty paragraphs right after the anchor paragraph:
    #   1) blank separator line
    #   2) "2 messages with length 64 and channel width of 100"
    #   3) first fragment data line
    #   4) second fragment data line
    $anchor.InsertParagraphAfter()
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()
    $anchor.Collapse(0)
    $anchor.InsertParagraphAfter()
    $anchor.Collapse(0)

    # Figure out where these four new (currently empty) paragraphs live.
    $p1 = $anchor.Paragraphs.Item(1)
    $blankPara = $d.Paragraphs.Item($p1.Index)
    $headingPara = $d.Paragraphs.Item($p1.Index + 1)
    $dataPara1 = $d.Paragraphs.Item($p1.Index + 2)
    $dataPara2 = $d.Paragraphs.Item($p1.Index + 3)

    # Leave the first paragraph (blank separator) empty; its indentation is
    # already inherited from the paragraph it was split from.

    # Heading line: "2 messages with length 64 and channel width of 100"
    $headingPara.Range.Text = "2 messages with length 64 and channel width of 100"

    # Data lines. Message length = 64, channel width = 100:
    #   8-byte header, 64 FF payload bytes, 4-byte "AA BB CC DD" footer,
    #   and (100 - 8 - 64 - 4) = 24 zero padding bytes.
    $ffPart = RepeatToken "FF" 64
    $padPart = RepeatToken "00" 24
    $footer = "AA BB CC DD"

    $prefix1 = "1, 1, "
    $header1 = "07 07 07 07 07 07 00 D5"
    $line1 = "$prefix1$header1 $ffPart $footer $padPart"

    $prefix2 = "2, 1, "
    $header2 = "07 07 07 07 07 07 01 D5"
    $line2 = "$prefix2$header2 $ffPart $footer $padPart"

    $dataPara1.Range.Text = $line1
    $dataPara2.Range.Text = $line2

    # Re-fetch the paragraphs (ranges are re-anchored after .Text assignment)
    # and colour the header (dark red) and footer (blue) runs, matching the
    # colour-coding used by the other examples in this section.
    $red = 192          # RGB(192,0,0)  == "C00000"
    $blue = 12419407    # RGB(79,129,189) == "4F81BD" (accent1)

    foreach ($pair in @(
            @{ Para = $d.Paragraphs.Item($p1.Index + 2); Prefix = $prefix1; Header = $header1 },
            @{ Para = $d.Paragraphs.Item($p1.Index + 3); Prefix = $prefix2; Header = $header2 }
        )) {
        $para = $pair.Para
        $prefix = $pair.Prefix
        $header = $pair.Header
        $pStart = $para.Range.Start

        $hdrStart = $pStart + $prefix.Length
        $hdrEnd = $hdrStart + $header.Length
        $hdrRange = $d.Range($hdrStart, $hdrEnd)
        $hdrRange.Font.Color = $red

        $ftrStart = $hdrEnd + 1 + $ffPart.Length + 1
        $ftrEnd = $ftrStart + $footer.Length
        $ftrRange = $d.Range($ftrStart, $ftrEnd)
        $ftrRange.Font.Color = $blue
    }
}
